$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 580; this shifts existing rows 580:633
# down to 581:634 (dimension grows from A1:T633 to A1:T634), matching the
# diff where every record previously at row N (N >= 580) now lives at
# row N+1, and a brand-new record appears at row 580.
$ws.Rows.Item(580).Insert()

# Populate the new row 580 with the new record's data.
$ws.Cells.Item(580, 1).Value = 6
$ws.Cells.Item(580, 2).Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Cells.Item(580, 3).Value = 'Metropolitana'
$ws.Cells.Item(580, 4).Value = 45021
$ws.Cells.Item(580, 5).Value = 13
$ws.Cells.Item(580, 6).Value = 'Fruta'
$ws.Cells.Item(580, 7).Value = 100101
$ws.Cells.Item(580, 8).Value = 'Berries'
$ws.Cells.Item(580, 9).Value = 100101001
$ws.Cells.Item(580, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(580, 11).Value = 'Sin especificar'
$ws.Cells.Item(580, 12).Value = 'Primera'
$ws.Cells.Item(580, 13).Value = 150
$ws.Cells.Item(580, 14).Value = 6000
$ws.Cells.Item(580, 15).Value = 6500
$ws.Cells.Item(580, 16).Value = 6250
$ws.Cells.Item(580, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(580, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(580, 19).Value = 3125
$ws.Cells.Item(580, 20).Value = 2
